$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Hydrogen - update B3 value, clear D3 (becomes an empty cell, no value)
$ws.Range("B3").Value = 3006369.219678119
$ws.Range("D3").ClearContents()
$ws.Range("D3").Style = "Normal"

# Row 4: Methanol - update C4 value
$ws.Range("C4").Value = 29.89584818596417

# Row 5: Ammonia - update C5 value
$ws.Range("C5").Value = 12117.66756755411

# Row 7: rename "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 4766.345533117365

# Row 8 (new): "Other" row, copying row 7's label formatting, plus D8 value
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 2899.317824044244
